$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 88888881
$ws.Range("A2").Value = 9999999
$ws.Range("A3").Value = 77777771
$ws.Range("A4").Value = 6666661
$ws.Range("A5").Value = 55555551

$ws.Range("A5").Select()
